# DOMA-6936: add "Is verified" column for contacts import example sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. New header cell H1 ---
$ws.Range("H1").Value = "Is verified"

# --- 2. New "Is verified" values for the 5 sample rows ---
$ws.Range("H2").Value = "Yes"
$ws.Range("H3").Value = "No"
$ws.Range("H4").Value = ""
$ws.Range("H5").Value = "yes"
$ws.Range("H6").Value = "no"

# --- 3. Fix the typo'd email in row 5 (was "test@example.com", now "ttest@example.com") ---
$ws.Range("F5").Value = "ttest@example.com"

# Rebuild hyperlinks so F5's mailto: link shows the corrected display text
# (the underlying mailto target stays test@example.com for every row).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("F2"), "mailto:test@example.com", "", "", "test@example.com")
$ws.Hyperlinks.Add($ws.Range("F3"), "mailto:test@example.com", "", "", "test1@example.com")
$ws.Hyperlinks.Add($ws.Range("F4"), "mailto:test@example.com", "", "", "test2@example.com")
$ws.Hyperlinks.Add($ws.Range("F5"), "mailto:test@example.com", "", "", "ttest@example.com")
$ws.Hyperlinks.Add($ws.Range("F6"), "mailto:test@example.com", "", "", "test3@example.com")

# --- 4. Match column H's width/style to the existing F/G columns ---
$refWidth = $ws.Columns.Item(6).ColumnWidth
$ws.Columns.Item(8).ColumnWidth = $refWidth

for ($r = 1; $r -le 6; $r++) {
    $ws.Cells.Item($r, 8).Style = $ws.Cells.Item($r, 7).Style
}
for ($r = 7; $r -le 50; $r++) {
    $ws.Cells.Item($r, 8).Style = $ws.Cells.Item($r, 6).Style
}
